$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 232
$ws1.Range("F3").Value = 55061
$ws1.Range("F4").Value = 1285
$ws1.Range("F6").Value = 877
$ws1.Range("F8").Value = 1166
$ws1.Range("F9").Value = 1447
$ws1.Range("F10").Value = 143
$ws1.Range("F11").Value = 48
$ws1.Range("F12").Value = 251
$ws1.Range("F13").Value = 417
$ws1.Range("F14").Value = 83
$ws1.Range("F15").Value = 44
$ws1.Range("F17").Value = 80
$ws1.Range("F19").Value = 5737
$ws1.Range("F21").Value = 5565
$ws1.Range("F22").Value = 9598
$ws1.Range("F24").Value = 166
$ws1.Range("F25").Value = 165
$ws1.Range("F26").Value = 257
$ws1.Range("F27").Value = 473
$ws1.Range("F28").Value = 152
$ws1.Range("F29").Value = 123
$ws1.Range("F30").Value = 4324
$ws1.Range("F31").Value = 338

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 60
$ws2.Range("F8").Value = 1156

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 232
$ws4.Range("F5").Value = 1285
$ws4.Range("F7").Value = 60
$ws4.Range("F8").Value = 877
$ws4.Range("F10").Value = 1166
$ws4.Range("F12").Value = 1447
$ws4.Range("F14").Value = 143
$ws4.Range("F15").Value = 251
$ws4.Range("F17").Value = 417
$ws4.Range("F18").Value = 83
$ws4.Range("F19").Value = 44
$ws4.Range("F22").Value = 80
$ws4.Range("F24").Value = 5737
$ws4.Range("F26").Value = 5565
$ws4.Range("F27").Value = 9598
$ws4.Range("F30").Value = 166
$ws4.Range("F31").Value = 165
$ws4.Range("F32").Value = 257
$ws4.Range("F34").Value = 473
$ws4.Range("F37").Value = 152
$ws4.Range("F38").Value = 123
$ws4.Range("F39").Value = 4324
$ws4.Range("F46").Value = 338
